# Update Name of Algo
# Applies updated RandomForest imputation results to the corresponding cells
# on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.34740000000001
$ws.Range("D7").Value = -7.106400000000002
$ws.Range("C10").Value = -13.04559999999999
$ws.Range("C12").Value = -10.71159999999999
$ws.Range("D15").Value = -7.992699999999996
$ws.Range("C18").Value = -11.47469999999999
$ws.Range("E18").Value = 18.18410000000001
$ws.Range("E19").Value = 16.56050000000001
$ws.Range("D20").Value = -7.859199999999998
$ws.Range("E27").Value = 16.38
$ws.Range("D29").Value = -6.831
$ws.Range("D30").Value = -7.7012
$ws.Range("D31").Value = -8.411400000000004
$ws.Range("C37").Value = -13.0906
$ws.Range("D40").Value = -8.413199999999993
$ws.Range("E42").Value = 16.6572
$ws.Range("E44").Value = 16.48229999999999
$ws.Range("E47").Value = 16.3381
$ws.Range("C55").Value = -13.6769
$ws.Range("E58").Value = 16.45280000000001
$ws.Range("C68").Value = -10.9582
$ws.Range("D68").Value = -7.024799999999995
$ws.Range("E73").Value = 17.43760000000001
$ws.Range("D76").Value = -7.604699999999999
$ws.Range("C77").Value = -12.1475
$ws.Range("C78").Value = -12.48380000000001
$ws.Range("D87").Value = -7.932899999999995
$ws.Range("D88").Value = -7.318899999999996
$ws.Range("E95").Value = 18.00780000000002
$ws.Range("D96").Value = -7.718500000000001
$ws.Range("D98").Value = -8.565700000000003
$ws.Range("D101").Value = -7.879700000000001
$ws.Range("E101").Value = 16.7815
$ws.Range("D102").Value = -8.047499999999998
